$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new Time Log entry as row 19, matching the formatting of the
# previous row (date format in column A, time format in columns B/C).
$ws.Range("A18:D18").Copy() | Out-Null
$ws.Range("A19:D19").PasteSpecial(-4122) | Out-Null

$ws.Range("A19").Value = 45792
$ws.Range("B19").Value = 0.20833333333333334
$ws.Range("C19").Value = 0.25
$ws.Range("D19").Value = "worked on overleaf and presentation"

$ws.Range("D19").Select() | Out-Null
